# Weekly refresh of the "Zapallo" price series: a new observation is
# inserted as row 90 (date 2022-04-29), which pushes every existing
# observation from row 90 onward down by one row (old row 90 -> new row
# 91, ..., old row 154 -> new row 155). The sheet's used range grows from
# A1:R154 to A1:R155 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 90; Excel shifts rows 90..154 down to 91..155 and
# copies the formatting (incl. the date number format on column D) from
# the row above, same as doing it interactively in the UI.
$ws.Rows.Item(90).Insert()

# Seed the new row from the record that is now directly below it (the
# record that used to be row 90), since most fields repeat verbatim for
# this market/product/variety/quality combination.
$ws.Range("A91:R91").Copy($ws.Range("A90:R90"))

# Overwrite the fields that actually differ for the new observation:
# a later date plus the corresponding volume/price figures.
$ws.Range("D90").Value = "2022-04-29"
$ws.Range("J90").Value = 160
$ws.Range("K90").Value = 300
$ws.Range("L90").Value = 350
$ws.Range("M90").Value = 325
$ws.Range("P90").Value = 325
